$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 32
$ws.Range("H32").Value = 2649.8572
$ws.Range("J32").Value = 3255.3333
$ws.Range("L32").Value = 3255.3333
$ws.Range("N32").Value = -3907.3333

# Row 51
$ws.Range("H51").Value = 6285.25
$ws.Range("I51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("M51").ClearContents()

# Row 55
$ws.Range("H55").Value = 1695.9231
$ws.Range("I55").Value = 166.83333
$ws.Range("K55").Value = 166.83333
$ws.Range("M55").Value = 47.16667000000001

# Row 70
$ws.Range("H70").Value = 2059.7334
$ws.Range("I70").Value = 1399.625
$ws.Range("J70").Value = 2814.1428
$ws.Range("K70").Value = 4198.875
$ws.Range("L70").Value = 8442.428400000001
$ws.Range("M70").Value = -3928.875
$ws.Range("N70").Value = -8982.428400000001

# Row 73
$ws.Range("H73").Value = 2059.7334
$ws.Range("I73").Value = 1399.625
$ws.Range("J73").Value = 2814.1428
$ws.Range("K73").Value = 4198.875
$ws.Range("L73").Value = 8442.428400000001
$ws.Range("M73").Value = -3262.875
$ws.Range("N73").Value = -10314.4284

# Row 112
$ws.Range("H112").Value = 5845.8184
$ws.Range("J112").Value = 6067.048
$ws.Range("L112").Value = 18201.144
$ws.Range("N112").Value = -20417.144

# Row 113
$ws.Range("H113").Value = 5325.533
$ws.Range("I113").Value = 2712.25
$ws.Range("J113").Value = 6275.8184
$ws.Range("K113").Value = 2712.25
$ws.Range("L113").Value = 6275.8184
$ws.Range("M113").Value = 541.75
$ws.Range("N113").Value = -12783.8184

# Row 116
$ws.Range("H116").Value = 5155.6924
$ws.Range("I116").Value = 4297.8
$ws.Range("J116").Value = 5691.875
$ws.Range("K116").Value = 4297.8
$ws.Range("L116").Value = 5691.875
$ws.Range("M116").Value = -855.8000000000002
$ws.Range("N116").Value = -12575.875

# Row 132
$ws.Range("H132").Value = 20410448
$ws.Range("I132").Value = 21741476
$ws.Range("K132").Value = 65224428
$ws.Range("M132").Value = -65221898

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 3380.7805
$ws.Range("I32").Value = 2761.4365
$ws.Range("K32").Value = 2761.4365
$ws.Range("M32").Value = -2474.4365

# Row 61
$ws.Range("H61").Value = 1863.7273
$ws.Range("I61").Value = 1084.6428
$ws.Range("J61").Value = 3227.125
$ws.Range("K61").Value = 1084.6428
$ws.Range("L61").Value = 3227.125
$ws.Range("M61").Value = -872.6428000000001
$ws.Range("N61").Value = -3651.125

# Row 122
$ws.Range("H122").Value = 510002.06
$ws.Range("I122").Value = 1940.697
$ws.Range("J122").Value = 2605755.2
$ws.Range("K122").Value = 5822.090999999999
$ws.Range("L122").Value = 7817265.600000001
$ws.Range("M122").Value = -3372.090999999999
$ws.Range("N122").Value = -7822165.600000001

# Row 136
$ws.Range("H136").Value = 1863.7273
$ws.Range("I136").Value = 1084.6428
$ws.Range("J136").Value = 3227.125
$ws.Range("K136").Value = 3253.9284
$ws.Range("L136").Value = 9681.375
$ws.Range("M136").Value = -703.9284000000002
$ws.Range("N136").Value = -14781.375

$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 238.07692
$ws.Range("I22").Value = 198
$ws.Range("J22").Value = 458.5
$ws.Range("K22").Value = 198
$ws.Range("L22").Value = 458.5
$ws.Range("M22").Value = 152
$ws.Range("N22").Value = -1158.5

# Row 86
$ws.Range("H86").Value = 7248.294
$ws.Range("I86").Value = 6409.8887
$ws.Range("K86").Value = 6409.8887
$ws.Range("M86").Value = -5286.8887

# Row 89
$ws.Range("H89").Value = 7248.294
$ws.Range("I89").Value = 6409.8887
$ws.Range("K89").Value = 32049.4435
$ws.Range("M89").Value = -26433.4435

# Row 99
$ws.Range("H99").Value = 3943.3333
$ws.Range("I99").Value = 2833.3333
$ws.Range("K99").Value = 2833.3333
$ws.Range("M99").Value = -1335.3333

# Row 122
$ws.Range("H122").Value = 3277
$ws.Range("I122").Value = 2836
$ws.Range("K122").Value = 8508
$ws.Range("M122").Value = -6058

# Row 126
$ws.Range("H126").Value = 3943.3333
$ws.Range("I126").Value = 2833.3333
$ws.Range("K126").Value = 8499.999899999999
$ws.Range("M126").Value = -6029.999899999999

# Row 134
$ws.Range("H134").Value = 3995.4644
$ws.Range("I134").Value = 3489
$ws.Range("K134").Value = 10467
$ws.Range("M134").Value = -7932

$ws = $wb.Worksheets.Item("CUL")
# Row 56
$ws.Range("H56").Value = 16672833
$ws.Range("I56").Value = 16672833
$ws.Range("K56").Value = 16672833
$ws.Range("M56").Value = -16672303

# Row 137
$ws.Range("H137").Value = 5178.375
$ws.Range("I137").Value = 1400
$ws.Range("J137").Value = 5718.143
$ws.Range("K137").Value = 4200
$ws.Range("L137").Value = 17154.429
$ws.Range("M137").Value = 900
$ws.Range("N137").Value = -27354.429

$ws = $wb.Worksheets.Item("GSM")
# Row 5
$ws.Range("H5").Value = 93.75
$ws.Range("I5").Value = 93.75
$ws.Range("K5").Value = 93.75
$ws.Range("M5").Value = 18.25

# Row 102
$ws.Range("H102").Value = 5439038
$ws.Range("I102").Value = 6537321
$ws.Range("K102").Value = 6537321
$ws.Range("M102").Value = -6535699

# Row 113
$ws.Range("H113").Value = 83334450
$ws.Range("J113").Value = 2222
$ws.Range("L113").Value = 2222
$ws.Range("N113").Value = -6562

# Row 122
$ws.Range("H122").Value = 400249.56
$ws.Range("I122").Value = 634194.2
$ws.Range("K122").Value = 1902582.6
$ws.Range("M122").Value = -1900132.6

# Row 132
$ws.Range("H132").Value = 3028.5144
$ws.Range("I132").Value = 2554.8215
$ws.Range("J132").Value = 4923.2856
$ws.Range("K132").Value = 7664.4645
$ws.Range("L132").Value = 14769.8568
$ws.Range("M132").Value = -5134.4645
$ws.Range("N132").Value = -19829.8568

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 296962.66
$ws.Range("I22").Value = 296962.66
$ws.Range("K22").Value = 296962.66
$ws.Range("M22").Value = -296667.66

# Row 27
$ws.Range("H27").Value = 296962.66
$ws.Range("I27").Value = 296962.66
$ws.Range("K27").Value = 296962.66
$ws.Range("M27").Value = -296855.66

# Row 48
$ws.Range("H48").Value = 31054
$ws.Range("I48").Value = 25135
$ws.Range("K48").Value = 25135
$ws.Range("M48").Value = -24474

$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 11906274
$ws.Range("I81").Value = 20834810
$ws.Range("K81").Value = 41669620
$ws.Range("M81").Value = -41668559

# Row 84
$ws.Range("H84").Value = 11906274
$ws.Range("I84").Value = 20834810
$ws.Range("K84").Value = 208348100
$ws.Range("M84").Value = -208342796

# Row 115
$ws.Range("H115").Value = 33688.5
$ws.Range("J115").Value = 33377
$ws.Range("L115").Value = 33377
$ws.Range("N115").Value = -36511

# Row 135
$ws.Range("H135").Value = 123993.25
$ws.Range("J135").Value = 132658
$ws.Range("L135").Value = 132658
$ws.Range("N135").Value = -142798

# Row 136
$ws.Range("H136").Value = 1984.8695
$ws.Range("I136").Value = 1355.8948
$ws.Range("K136").Value = 4067.6844
$ws.Range("M136").Value = -1517.6844
